$d = $word.ActiveDocument
$d.Content.Find.Execute("Sept", $true, $false, $false, $false, $false, $true, 1, $false, "Oct", 2)
